$d = $word.ActiveDocument

# The commit appends seven blank "Normal" style paragraphs (justified,
# 1.5-ish line spacing = 360 twips "auto") right after the document's
# final paragraph (the one holding the _gjdgxs bookmark), before sectPr.
$paraXml = "<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>" +
           "<w:pPr>" +
           "<w:pStyle w:val='Normal'/>" +
           "<w:spacing w:line='360' w:lineRule='auto'/>" +
           "<w:jc w:val='both'/>" +
           "<w:rPr><w:lang w:val='en-US'/></w:rPr>" +
           "</w:pPr>" +
           "</w:p>"

for ($i = 0; $i -lt 7; $i++) {
    $endPos = $d.Paragraphs.Last.Range.End
    $insertPoint = $d.Range($endPos, $endPos)
    $insertPoint.InsertXML($paraXml)
}

Write-Output $d.Paragraphs.Count
